$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.74%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.72%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.149"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.23%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07368"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.31%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.429"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "59.03%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.933"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.17%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.767"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.68%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9191"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.05%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1727"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.55%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07508"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.82%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08119"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.73%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03031"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.84%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09919"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.21%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001492"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.53%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006079"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-6.25%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.460"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.08%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.229"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.05%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3278"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.42%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.10%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.652"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.71%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04642"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.89%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1568"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.24%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.59%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004477"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.74%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.07%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "10.71%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01731"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.31%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04523"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.62%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007224"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.68%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1342"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.41%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002231"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.77%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01085"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-15.01%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006285"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.87%"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.695"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "138.83%"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01000"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-23.02%"
